$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before B (splits "name" into first_name/last_name) ---
$ws.Columns("B").Insert()

# --- 2. Remove the stale hyperlinks left over from the old F-column emails ---
#     (Insert() does not re-anchor hyperlink ranges, so start clean.)
$links = @($ws.Hyperlinks)
for ($i = $links.Count - 1; $i -ge 0; $i--) {
  $links[$i].Delete()
}

# --- 3. Header row ---
$ws.Range("A1").Value = "first_name"
$ws.Range("B1").Value = "last_name"
$ws.Range("C1").Value = "mssv"
$ws.Range("D1").Value = "phone"
$ws.Range("E1").Value = "birthday"
$ws.Range("F1").Value = "address"
$ws.Range("G1").Value = "email"
$ws.Range("H1").Value = "sex"

# --- 4. Row 2 data (replaces the old Hai/Nam/Tho sample rows) ---
$ws.Range("A2").Value = "Trần Văn"
$ws.Range("B2").Value = "Khang"
$ws.Range("C2").Value = 20140956
$ws.Range("D2").Value = 987884444
$ws.Range("E2").Value = 35371
$ws.Range("F2").Value = "Hà Giang"
$ws.Range("G2").Value = "khangtv20140956@gmail.com"
$ws.Range("H2").Value = "nam"

# --- 5. Clear out rows 3 & 4 entirely (keep their row heights) ---
$ws.Range("A3:H4").Clear()
$ws.Rows(3).RowHeight = 75
$ws.Rows(4).RowHeight = 75
$ws.Rows(2).RowHeight = 31.5

# --- 6. Re-create the single hyperlink for the new e-mail address, then strip
#        the auto-applied hyperlink look (the source keeps plain formatting) ---
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:khangtv20140956@gmail.com") | Out-Null
$ws.Range("G2").Style = "Normal"

# --- 7. Column widths (character widths; Excel rounds to the nearest pixel,
#        these are the closest achievable inputs for the target widths) ---
$ws.Columns(1).ColumnWidth = 12.6666666666667
$ws.Columns(2).ColumnWidth = 13.6666666666667
$ws.Columns(3).ColumnWidth = 12.3333333333333
$ws.Columns(4).ColumnWidth = 12.6666666666667
$ws.Columns(5).ColumnWidth = 8.83333333333333
$ws.Columns(6).ColumnWidth = 15.8333333333333
$ws.Columns(7).ColumnWidth = 30
$ws.Columns(9).ColumnWidth = 16

# --- 8. Touch I6 (format-only) so the sheet's used range / dimension grows
#        to A1:I6, matching the widened, taller layout (column I carries a
#        width but no header, same as column H did before the edit). ---
$ws.Range("I6").Style = "Normal"
